$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Baru")

# --- Restructure the "Baru" sheet ---
# Original layout:
#   row1: namaMobilBaru | condition               (header)
#   row2: Toyota Rush   | passed                  (s=2)
#   row3: Daihatsu Sigra| passed                  (s=2)
#   row4: Honda Brio    | failed                  (s=1)
#
# Target layout:
#   row1: namaMobilBaru | condition | spesific               | paket   (header, C/D centered s=6)
#   row2: Rush          | passed    | Rush 1.5 G M/T          | yes     (s=2)
#   row3: Rush          | passed    | Rush 1.5 S TRD A/T      | no      (s=2)
#   row4: Rush          | passed    | Rush 1.5 G M/T          | all     (s=2)
#   row5: Sigra         | passed    | Sigra 1.2 R A/T MC      | yes     (s=2)
#   row6: Sigra         | passed    | Sigra 1.0 D M/T MC      | no      (s=2)
#   row7: Sigra         | passed    | Sigra 1.2 R M/T DLX MC  | all     (s=2)
#   row8: Brio          | failed    |                         |         (s=1)

# Insert two extra rows after the "Rush" row (row2) so the Rush family has 3 rows (2-4).
$ws.Rows.Item(3).Insert() | Out-Null
$ws.Rows.Item(3).Insert() | Out-Null

# The former "Daihatsu Sigra" row is now row5. Insert two extra rows after it for the Sigra family (5-7).
$ws.Rows.Item(6).Insert() | Out-Null
$ws.Rows.Item(6).Insert() | Out-Null

# --- Fill in the new data (order chosen to mirror the author's natural entry sequence) ---
$ws.Range("D1").Value = "paket"
$ws.Range("D2").Value = "yes"
$ws.Range("A8").Value = "Brio"
$ws.Range("C1").Value = "spesific"
$ws.Range("C2").Value = "Rush 1.5 G M/T"
$ws.Range("D3").Value = "no"
$ws.Range("D4").Value = "all"
$ws.Range("C3").Value = "Rush 1.5 S TRD A/T"
$ws.Range("C7").Value = "Sigra 1.2 R M/T DLX MC"
$ws.Range("C6").Value = "Sigra 1.0 D M/T MC"
$ws.Range("C5").Value = "Sigra 1.2 R A/T MC"

$ws.Range("C4").Value = "Rush 1.5 G M/T"
$ws.Range("D5").Value = "yes"
$ws.Range("D6").Value = "no"
$ws.Range("D7").Value = "all"

$ws.Range("A2").Value = "Rush"
$ws.Range("A3").Value = "Rush"
$ws.Range("A4").Value = "Rush"
$ws.Range("A5").Value = "Sigra"
$ws.Range("A6").Value = "Sigra"
$ws.Range("A7").Value = "Sigra"

$ws.Range("B2").Value = "passed"
$ws.Range("B3").Value = "passed"
$ws.Range("B4").Value = "passed"
$ws.Range("B5").Value = "passed"
$ws.Range("B6").Value = "passed"
$ws.Range("B7").Value = "passed"

# Apply the same row fill (s=2 / "passed" green) to the new C:D cells for rows 2-7.
$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("C2:D7").PasteSpecial(-4122) | Out-Null
$ws.Range("C3:D3").PasteSpecial(-4122) | Out-Null

# Row 8 ("Brio"/"failed") keeps style s=1 on C8:D8 (blank, styled like the rest of the row).
$ws.Range("A8:B8").Copy() | Out-Null
$ws.Range("C8:D8").PasteSpecial(-4122) | Out-Null

# Header row: center-align the new C1/D1 header cells (creates the new centered style).
$ws.Range("C1:D1").HorizontalAlignment = -4108

$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null

$ws.Application.CutCopyMode = $false

$ws.Range("F8").Select() | Out-Null
